$d = $word.ActiveDocument

# 1. Merge the "are" run (with its proofErr gramStart/gramEnd wrapper) and the
#    following run back into a single run, matching the collapsed text.
$d.Content.Find.Execute(" are important which", $true, $false, $false, $false, $false, $true, 1, $false, " are important which", 2)

# 2. Remove the whole "With the below model highest auc score on leaderboard..." bullet paragraph.
$p = $d.Paragraphs(24)
$r = $p.Range
$r.Delete()

# 3. Re-home the "_GoBack" bookmark onto the now-empty paragraph that used to follow
#    the deleted bullet (instead of after the "Contact No." paragraph).
$bms = $d.Bookmarks
$old = $bms.Item("_GoBack")
$old.Delete()

$target = $d.Paragraphs(24).Range
$d.Bookmarks.Add("_GoBack", $target)
